# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Swap Alemania (row 24) / Turquia (row 25) order: row 24 becomes Turquia, row 25 becomes Alemania
$ws.Range("A24").Value = "Turquia"
$ws.Range("A25").Value = "Alemania"

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Octubre de 2020 a las 20:11"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 8110587
$ws.Range("C4").Value = 20334
$ws.Range("D4").Value = 5243662
$ws.Range("E4").Value = 2645691
$ws.Range("G4").Value = 361
$ws.Range("H4").Value = 221234

# Row 24 - now Turquia
$ws.Range("B24").Value = 340450
$ws.Range("C24").Value = 1671
$ws.Range("D24").Value = 298368
$ws.Range("E24").Value = 33068
$ws.Range("G24").Value = 57
$ws.Range("H24").Value = 9014

# Row 25 - now Alemania
$ws.Range("B25").Value = 339722
$ws.Range("C25").Value = 4043
$ws.Range("D25").Value = 281900
$ws.Range("E25").Value = 48063
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = 9759

# Row 33 - Marruecos
$ws.Range("B33").Value = 160333
$ws.Range("C33").Value = 3387
$ws.Range("D33").Value = 133959
$ws.Range("E33").Value = 23648
$ws.Range("G33").Value = 41
$ws.Range("H33").Value = 2726

# Row 65 - Libano
$ws.Range("B65").Value = 57246
$ws.Range("C65").Value = 1377
$ws.Range("D65").Value = 25164
$ws.Range("E65").Value = 31583
$ws.Range("G65").Value = 20
$ws.Range("H65").Value = 499

# Row 66 - Argelia
$ws.Range("B66").Value = 53584
$ws.Range("C66").Value = 185
$ws.Range("D66").Value = 37603
$ws.Range("E66").Value = 14154
$ws.Range("G66").Value = 9
$ws.Range("H66").Value = 1827

# Row 72 - Irlanda
$ws.Range("B72").Value = 45243
$ws.Range("C72").Value = 1084
$ws.Range("E72").Value = 20044
$ws.Range("G72").Value = 5
$ws.Range("H72").Value = 1835

# Row 105 - Maldivas
$ws.Range("B105").Value = 11062
$ws.Range("C105").Value = 69
$ws.Range("D105").Value = 9880
$ws.Range("E105").Value = 1147

# Row 128 - Sri Lanka
$ws.Range("B128").Value = 5170
$ws.Range("C128").Value = 132
$ws.Range("E128").Value = 1800

# Row 150 - Principado de Andorra
$ws.Range("B150").Value = 3190
$ws.Range("C150").Value = 195
$ws.Range("D150").Value = 2011
$ws.Range("E150").Value = 1120
$ws.Range("G150").Value = 2
$ws.Range("H150").Value = 59

# Row 189 - Monaco
$ws.Range("B189").Value = 248
$ws.Range("C189").Value = 7
$ws.Range("E189").Value = 29
